$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Content edits on "Arduino Pins" sheet (performed BEFORE the row sort,
#    so that the sort picks up the corrected values).
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Arduino Pins")

# Pin 0/TX0 should really be pin 1 (Arduino "1/TX0")
$ws2.Range("D73").Value = "1/TX0"

# Note that pin 4 is driven by both PA29 and PC26
$ws2.Range("E71").Value = "SD-CS, wired to PC26"
$ws2.Range("E72").Value = "Wired to PA29"

# The four "Native USB" rows (DFSDM/DFSDP/DHSDM/DHSDP) had their value in the
# wrong column (Peripheral B instead of Extra Function) - move it over.
$ws2.Range("D77").Value = "Native USB"
$ws2.Range("C77").ClearContents()
$ws2.Range("D78").Value = "Native USB"
$ws2.Range("C78").ClearContents()
$ws2.Range("D79").Value = "Native USB"
$ws2.Range("C79").ClearContents()
$ws2.Range("D80").Value = "Native USB"
$ws2.Range("C80").ClearContents()

# ---------------------------------------------------------------------------
# 2. Sort the data rows (2-80) by column D (Arduino Pin) ascending.
# ---------------------------------------------------------------------------
$dataRange = $ws2.Range("A2:E80")
$sortKey = $ws2.Range("D2:D80")
$dataRange.Sort($sortKey, 1)

# ---------------------------------------------------------------------------
# 3. View/pane/selection updates on "Arduino Pins".
# ---------------------------------------------------------------------------
$ws2.Application.ActiveWindow.FreezePanes = $false
$ws2.Range("A25").Select()
$ws2.Application.ActiveWindow.FreezePanes = $true
$ws2.Range("E39").Select()

# ---------------------------------------------------------------------------
# 4. Column width updates.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Chip Pin Mapping")
$ws1.Range("A1:B1").EntireColumn.ColumnWidth = 19
$ws1.Range("C1").EntireColumn.ColumnWidth = 45.0408163265306
$ws1.Range("D1:G1").EntireColumn.ColumnWidth = 19
$ws1.Range("H1:J1").EntireColumn.ColumnWidth = 22.6530612244898
$ws1.Range("K1:IV1").EntireColumn.ColumnWidth = 19

$ws2.Range("A1:B1").EntireColumn.ColumnWidth = 19
$ws2.Range("C1:E1").EntireColumn.ColumnWidth = 22.6530612244898
$ws2.Range("F1:IQ1").EntireColumn.ColumnWidth = 19
